$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2,1).Value = "Última actualización: 08:31:16"
$ws.Cells.Item(3,1).Value = "Total filas: 83"
$ws.Cells.Item(35,3).Value = "16_SANTA ANA"
$ws.Cells.Item(36,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(38,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(39,3).Value = "27_EL RETIRO"
$ws.Cells.Item(44,3).Value = "16_SANTA ANA"
$ws.Cells.Item(45,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(55,1).Value = "08:31:16"
$ws.Cells.Item(55,2).Value = "08:32"
$ws.Cells.Item(55,3).Value = "10_OLMOS"
$ws.Cells.Item(55,4).Value = 1
$ws.Cells.Item(56,1).Value = "07:59:51"
$ws.Cells.Item(56,2).Value = "08:34"
$ws.Cells.Item(56,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(56,4).Value = 35
$ws.Cells.Item(57,1).Value = "08:31:16"
$ws.Cells.Item(57,2).Value = "08:35"
$ws.Cells.Item(57,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(57,4).Value = 4
$ws.Cells.Item(58,1).Value = "08:31:16"
$ws.Cells.Item(58,2).Value = "08:42"
$ws.Cells.Item(58,4).Value = 11
$ws.Cells.Item(59,2).Value = "08:43"
$ws.Cells.Item(59,3).Value = "14_ABASTO"
$ws.Cells.Item(59,4).Value = 44
$ws.Cells.Item(60,1).Value = "07:23:38"
$ws.Cells.Item(60,2).Value = "08:44"
$ws.Cells.Item(60,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(60,4).Value = 81
$ws.Cells.Item(61,1).Value = "08:31:16"
$ws.Cells.Item(61,2).Value = "08:44"
$ws.Cells.Item(61,3).Value = "14_ABASTO"
$ws.Cells.Item(61,4).Value = 13
$ws.Cells.Item(62,1).Value = "08:31:16"
$ws.Cells.Item(62,2).Value = "08:53"
$ws.Cells.Item(62,3).Value = "10_OLMOS"
$ws.Cells.Item(62,4).Value = 22
$ws.Cells.Item(63,1).Value = "08:31:16"
$ws.Cells.Item(63,2).Value = "08:54"
$ws.Cells.Item(63,3).Value = "17_ROMERO"
$ws.Cells.Item(63,4).Value = 23
$ws.Cells.Item(64,1).Value = "08:31:16"
$ws.Cells.Item(64,2).Value = "09:01"
$ws.Cells.Item(64,3).Value = "215A_EL PATO"
$ws.Cells.Item(64,4).Value = 30
$ws.Cells.Item(65,2).Value = "09:03"
$ws.Cells.Item(65,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(65,4).Value = 64
$ws.Cells.Item(66,1).Value = "08:31:16"
$ws.Cells.Item(66,2).Value = "09:04"
$ws.Cells.Item(66,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(66,4).Value = 33
$ws.Cells.Item(67,1).Value = "08:31:16"
$ws.Cells.Item(67,2).Value = "09:05"
$ws.Cells.Item(67,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(67,4).Value = 34
$ws.Cells.Item(68,2).Value = "09:10"
$ws.Cells.Item(68,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(68,4).Value = 71
$ws.Cells.Item(69,1).Value = "08:31:16"
$ws.Cells.Item(69,2).Value = "09:11"
$ws.Cells.Item(69,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(69,4).Value = 40
$ws.Cells.Item(70,2).Value = "09:16"
$ws.Cells.Item(70,3).Value = "27_EL RETIRO"
$ws.Cells.Item(70,4).Value = 77
$ws.Cells.Item(71,1).Value = "08:31:16"
$ws.Cells.Item(71,2).Value = "09:17"
$ws.Cells.Item(71,3).Value = "27_EL RETIRO"
$ws.Cells.Item(71,4).Value = 46
$ws.Cells.Item(72,2).Value = "09:20"
$ws.Cells.Item(72,3).Value = "81_EL PELIGRO"
$ws.Cells.Item(72,4).Value = 81
$ws.Cells.Item(73,1).Value = "08:31:16"
$ws.Cells.Item(73,2).Value = "09:21"
$ws.Cells.Item(73,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(73,4).Value = 50
$ws.Cells.Item(73,5).Value = "LP1912"
$ws.Cells.Item(74,1).Value = "07:59:51"
$ws.Cells.Item(74,2).Value = "09:22"
$ws.Cells.Item(74,3).Value = "17_ROMERO"
$ws.Cells.Item(74,4).Value = 83
$ws.Cells.Item(74,5).Value = "LP1912"
$ws.Cells.Item(75,1).Value = "08:31:16"
$ws.Cells.Item(75,2).Value = "09:23"
$ws.Cells.Item(75,3).Value = "16_SANTA ANA"
$ws.Cells.Item(75,4).Value = 52
$ws.Cells.Item(75,5).Value = "LP1912"
$ws.Cells.Item(76,1).Value = "07:59:51"
$ws.Cells.Item(76,2).Value = "09:23"
$ws.Cells.Item(76,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(76,4).Value = 84
$ws.Cells.Item(76,5).Value = "LP1912"
$ws.Cells.Item(77,1).Value = "08:31:16"
$ws.Cells.Item(77,2).Value = "09:24"
$ws.Cells.Item(77,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(77,4).Value = 53
$ws.Cells.Item(77,5).Value = "LP1912"
$ws.Cells.Item(78,1).Value = "08:31:16"
$ws.Cells.Item(78,2).Value = "09:32"
$ws.Cells.Item(78,3).Value = "15_ABASTO"
$ws.Cells.Item(78,4).Value = 61
$ws.Cells.Item(78,5).Value = "LP1912"
$ws.Cells.Item(79,1).Value = "08:31:16"
$ws.Cells.Item(79,2).Value = "09:33"
$ws.Cells.Item(79,3).Value = "10_OLMOS"
$ws.Cells.Item(79,4).Value = 62
$ws.Cells.Item(79,5).Value = "LP1912"
$ws.Cells.Item(80,1).Value = "08:31:16"
$ws.Cells.Item(80,2).Value = "09:35"
$ws.Cells.Item(80,3).Value = "16_SANTA ANA"
$ws.Cells.Item(80,4).Value = 64
$ws.Cells.Item(80,5).Value = "LP1912"
$ws.Cells.Item(81,1).Value = "07:59:51"
$ws.Cells.Item(81,2).Value = "09:41"
$ws.Cells.Item(81,3).Value = "215C_EL PATO"
$ws.Cells.Item(81,4).Value = 102
$ws.Cells.Item(81,5).Value = "LP1912"
$ws.Cells.Item(82,1).Value = "08:31:16"
$ws.Cells.Item(82,2).Value = "09:42"
$ws.Cells.Item(82,3).Value = "215C_EL PATO"
$ws.Cells.Item(82,4).Value = 71
$ws.Cells.Item(82,5).Value = "LP1912"
$ws.Cells.Item(83,1).Value = "07:59:51"
$ws.Cells.Item(83,2).Value = "09:43"
$ws.Cells.Item(83,3).Value = "14_ABASTO"
$ws.Cells.Item(83,4).Value = 104
$ws.Cells.Item(83,5).Value = "LP1912"
$ws.Cells.Item(84,1).Value = "08:31:16"
$ws.Cells.Item(84,2).Value = "09:44"
$ws.Cells.Item(84,3).Value = "14_ABASTO"
$ws.Cells.Item(84,4).Value = 73
$ws.Cells.Item(84,5).Value = "LP1912"
$ws.Cells.Item(85,1).Value = "08:31:16"
$ws.Cells.Item(85,2).Value = "10:12"
$ws.Cells.Item(85,3).Value = "15_ABASTO"
$ws.Cells.Item(85,4).Value = 101
$ws.Cells.Item(85,5).Value = "LP1912"
$ws.Cells.Item(86,1).Value = "08:31:16"
$ws.Cells.Item(86,2).Value = "10:21"
$ws.Cells.Item(86,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(86,4).Value = 110
$ws.Cells.Item(86,5).Value = "LP1912"
$ws.Cells.Item(87,1).Value = "08:31:16"
$ws.Cells.Item(87,2).Value = "10:22"
$ws.Cells.Item(87,3).Value = "17_ROMERO"
$ws.Cells.Item(87,4).Value = 111
$ws.Cells.Item(87,5).Value = "LP1912"
$ws.Cells.Item(88,1).Value = "08:31:16"
$ws.Cells.Item(88,2).Value = "10:26"
$ws.Cells.Item(88,3).Value = "215A_EL PATO"
$ws.Cells.Item(88,4).Value = 115
$ws.Cells.Item(88,5).Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2,1).Value = "Última actualización: 08:31:16"
$ws.Cells.Item(3,1).Value = "Total filas: 11"
$ws.Cells.Item(13,1).Value = "08:31:16"
$ws.Cells.Item(13,4).Value = 30
$ws.Cells.Item(15,1).Value = "08:31:16"
$ws.Cells.Item(15,2).Value = "09:42"
$ws.Cells.Item(15,3).Value = "215C_EL PATO"
$ws.Cells.Item(15,4).Value = 71
$ws.Cells.Item(15,5).Value = "LP1912"
$ws.Cells.Item(16,1).Value = "08:31:16"
$ws.Cells.Item(16,2).Value = "10:26"
$ws.Cells.Item(16,3).Value = "215A_EL PATO"
$ws.Cells.Item(16,4).Value = 115
$ws.Cells.Item(16,5).Value = "LP1912"

# --- Sheet 3: 6203-6173 ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2,1).Value = "Última actualización: 08:31:16"
$ws.Cells.Item(3,1).Value = "Total filas: 16"
$ws.Cells.Item(18,1).Value = "08:31:16"
$ws.Cells.Item(18,2).Value = "08:37"
$ws.Cells.Item(18,3).Value = "215A_LA PLATA"
$ws.Cells.Item(18,4).Value = 6
$ws.Cells.Item(18,5).Value = "L6173"
$ws.Cells.Item(19,1).Value = "07:59:51"
$ws.Cells.Item(19,2).Value = "09:08"
$ws.Cells.Item(19,3).Value = "215D_LA PLATA"
$ws.Cells.Item(19,4).Value = 69
$ws.Cells.Item(19,5).Value = "L6203"
$ws.Cells.Item(20,1).Value = "08:31:16"
$ws.Cells.Item(20,2).Value = "09:09"
$ws.Cells.Item(20,3).Value = "215D_LA PLATA"
$ws.Cells.Item(20,4).Value = 38
$ws.Cells.Item(20,5).Value = "L6203"
$ws.Cells.Item(21,1).Value = "08:31:16"
$ws.Cells.Item(21,2).Value = "10:03"
$ws.Cells.Item(21,3).Value = "215B_LP-P MOR-40 Y 115"
$ws.Cells.Item(21,4).Value = 92
$ws.Cells.Item(21,5).Value = "L6173"
